$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers for ci.lower and ci.upper columns
$ws.Range("G1").Value = "ci.lower"
$ws.Range("H1").Value = "ci.upper"

# Fill in the confidence interval values for rows 2-26
$ws.Range("G2").Value = -0.441440529068733
$ws.Range("H2").Value = 0.121593735710427
$ws.Range("G3").Value = -0.0106823450117037
$ws.Range("H3").Value = 0.0782141872810027
$ws.Range("G4").Value = 0.00251075376497285
$ws.Range("H4").Value = 0.108252652338831
$ws.Range("G5").Value = -0.0620540692612157
$ws.Range("H5").Value = 0.454348608795867
$ws.Range("G6").Value = 0.0127400165824099
$ws.Range("H6").Value = 0.549293445309827
$ws.Range("G7").Value = -0.0578640860207077
$ws.Range("H7").Value = 0.423670313579011
$ws.Range("G8").Value = 0.0136002414894887
$ws.Range("H8").Value = 0.586382557391756
$ws.Range("G9").Value = -0.000416615929511546
$ws.Range("H9").Value = 0.0895642401160629
$ws.Range("G10").Value = -0.46761987925916
$ws.Range("H10").Value = -0.391709231977247
$ws.Range("G11").Value = -0.0596665803999276
$ws.Range("H11").Value = 0.016435016565423
# Row 12: G/H left empty (matches blank C:F in source row)
# Row 13: G/H left empty (matches blank C:F in source row)
# Row 14: G/H left empty (matches blank C:F in source row)
$ws.Range("G15").Value = -0.0744647363850066
$ws.Range("H15").Value = 0.0285486760329598
$ws.Range("G16").Value = -0.403360302159362
$ws.Range("H16").Value = 0.154642360262529
$ws.Range("G17").Value = -0.461223735163437
$ws.Range("H17").Value = 0.153798656807579
$ws.Range("G18").Value = -0.156428374562005
$ws.Range("H18").Value = 0.0213646900234074
$ws.Range("G19").Value = -0.847340627158023
$ws.Range("H19").Value = 0.115728172041415
$ws.Range("G20").Value = -0.908697217591734
$ws.Range("H20").Value = 0.124108138522431
$ws.Range("G21").Value = -0.016435016565423
$ws.Range("H21").Value = 0.0596665803999276
$ws.Range("G22").Value = -0.000833231859023134
$ws.Range("H22").Value = 0.179128480232126
$ws.Range("G23").Value = -0.0890251355157945
$ws.Range("H23").Value = 0.323201706838736
$ws.Range("G24").Value = -0.00451344717970481
$ws.Range("H24").Value = 0.970302473619253
$ws.Range("G25").Value = -0.116904874184011
$ws.Range("H25").Value = 0.286643796541596
$ws.Range("G26").Value = -0.00931450514640708
$ws.Range("H26").Value = 0.963642506573295
